$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.838.15"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.214.57"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.74"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.615"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.90"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.02%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -2.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.63"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.21%  "
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("E12").Value = "  -1.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.547.11"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.20"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.10%  "
$ws.Range("E16").Value = "  -2.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.206.08"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.768.76"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("E19").Value = "  +9.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.96"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.11%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.28"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +17.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "229.16"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.92%  "
$ws.Range("E24").Value = "  -5.67%  "
$ws.Range("E25").Value = "  +1.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("E28").Value = "  -1.80%  "
$ws.Range("E29").Value = "  +0.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.80"
$ws.Range("D30").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.53"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.63"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +8.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0794"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.12%  "
$ws.Range("E34").Value = "  -0.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "29.17"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -6.79%  "
$ws.Range("E36").Value = "  -8.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.25"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.15%  "
$ws.Range("E38").Value = "  -4.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.38"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "65.81"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +5.97%  "
$ws.Range("E41").Value = "  -3.18%  "
$ws.Range("E42").Value = "  -2.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.197"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.65"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.65"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.61%  "
$ws.Range("E46").Value = "  -2.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.37"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.87%  "
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("E49").Value = "  -0.38%  "
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.420.31"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.38%  "
